# Update F/G columns (想去人数 / 最低票价) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 290
    $ws.Range("F3").Value = 298
    $ws.Range("F4").Value = 863
    $ws.Range("F6").Value = 311
    $ws.Range("F7").Value = 9038
    $ws.Range("G7").Value = 65
    $ws.Range("F8").Value = 78
    $ws.Range("F11").Value = 119
    $ws.Range("F17").Value = 267
    $ws.Range("F18").Value = 764
    $ws.Range("F19").Value = 39
    $ws.Range("F20").Value = 88
}
